$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data as of the Tue Oct  1 09:56:23 UTC 2024 GitHub Actions refresh.
# Columns: Row|~|Coin|~|Link|~|Price|~|Volume(1h)
$rowsData = @"
2|~|Bitcoin|~|https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc|~|63.955.26|~|  +0.26%  
3|~|Ethereum|~|https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth|~|2.639.50|~|  +0.89%  
4|~|TetherUSD|~|https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt|~|1.00|~|  -0.09%  
5|~|BNB|~|https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb|~|579.41|~|  +0.55%  
6|~|Solana|~|https://coinranking.com/coin/zNZHO_Sjf+solana-sol|~|156.99|~|  +0.91%  
7|~|XRP|~|https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp|~|0.630|~|  -1.85%  
8|~|USDC|~|https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc|~|1.00|~|  -0.07%  
9|~|Dogecoin|~|https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge|~|0.118|~|  -2.10%  
10|~|Toncoin|~|https://coinranking.com/coin/67YlI0K1b+toncoin-ton|~|5.84|~|  +0.98%  
11|~|Cardano|~|https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada|~|0.384|~|  -0.29%  
12|~|TRON|~|https://coinranking.com/coin/qUhEFk1I61atv+tron-trx|~|0.157|~|  +1.00%  
13|~|Avalanche|~|https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax|~|28.69|~|  +1.13%  
14|~|WrappedliquidstakedEther2.0|~|https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth|~|3.114.41|~|  +0.56%  
15|~|ShibaInu|~|https://coinranking.com/coin/xz24e0BjL+shibainu-shib|~|0.0000185|~|  +0.92%  
16|~|WrappedBTC|~|https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc|~|63.853.68|~|  +0.32%  
17|~|WrappedEther|~|https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth|~|2.642.92|~|  +0.56%  
18|~|Chainlink|~|https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link|~|12.16|~|  +0.29%  
19|~|Uniswap|~|https://coinranking.com/coin/_H5FVG9iW+uniswap-uni|~|7.76|~|  +3.20%  
20|~|Polkadot|~|https://coinranking.com/coin/25W7FG7om+polkadot-dot|~|4.53|~|  -2.29%  
21|~|BitcoinCash|~|https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch|~|345.25|~|  +0.28%  
22|~|Dai|~|https://coinranking.com/coin/MoTuySvg7+dai-dai|~|1.00|~|  +0.19%  
23|~|Litecoin|~|https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc|~|68.28|~|  +1.12%  
24|~|SuiNetwork|~|https://coinranking.com/coin/3xJluUMvp+suinetwork-sui|~|1.88|~|  +8.13%  
25|~|PEPE|~|https://coinranking.com/coin/03WI8NQPF+pepe-pepe|~|0.0000113|~|  +3.35%  
26|~|Fetch.AI|~|https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet|~|1.64|~|  +4.65%  
27|~|InternetComputer(DFINITY)|~|https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp|~|9.29|~|  -0.12%  
28|~|Bittensor|~|https://coinranking.com/coin/pgv7xSFi6+bittensor-tao|~|582.28|~|  +1.67%  
29|~|Aptos|~|https://coinranking.com/coin/HGYj5JCv5+aptos-apt|~|8.22|~|  +4.23%  
30|~|Kaspa|~|https://coinranking.com/coin/V8GxkwWow+kaspa-kas|~|0.161|~|  +0.55%  
31|~|Binance-PegBSC-USD|~|https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd|~|1.00|~|  -0.13%  
32|~|PancakeSwap|~|https://coinranking.com/coin/ncYFcP709+pancakeswap-cake|~|2.06|~|  -0.49%  
33|~|ImmutableX|~|https://coinranking.com/coin/Z96jIvLU7+immutablex-imx|~|1.74|~|  +2.16%  
34|~|RenderToken|~|https://coinranking.com/coin/vfo5XYwcV+rendertoken-render|~|6.63|~|  +2.97%  
35|~|NEARProtocol|~|https://coinranking.com/coin/DCrsaMv68+nearprotocol-near|~|5.46|~|  +3.03%  
36|~|PolygonEcosystemToken|~|https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol|~|0.404|~|  -1.17%  
37|~|EthereumClassic|~|https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc|~|19.79|~|  -0.61%  
38|~|FirstDigitalUSD|~|https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd|~|0.999|~|  +0.02%  
39|~|Stacks|~|https://coinranking.com/coin/mMPrMcB7+stacks-stx|~|1.92|~|  +2.53%  
40|~|Monero|~|https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr|~|152.90|~|  +0.97%  
41|~|dogwifhat|~|https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif|~|2.55|~|  +7.78%  
42|~|USDe|~|https://coinranking.com/coin/exbfr2U-0+usde-usde|~|0.999|~|  -0.04%  
43|~|Aave|~|https://coinranking.com/coin/ixgUfzmLR+aave-aave|~|162.42|~|  +4.34%  
44|~|InjectiveProtocol|~|https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj|~|24.23|~|  +5.06%  
45|~|Filecoin|~|https://coinranking.com/coin/ymQub4fuB+filecoin-fil|~|3.92|~|  -1.06%  
46|~|Hedera|~|https://coinranking.com/coin/jad286TjB+hedera-hbar|~|0.0590|~|  -0.67%  
47|~|Mantle|~|https://coinranking.com/coin/BoI4ux0nd+mantle-mnt|~|0.635|~|  +0.84%  
48|~|Stellar|~|https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm|~|0.100|~|  -1.60%  
49|~|VeChain|~|https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet|~|0.0249|~|  -0.67%  
50|~|EnergySwap|~|https://coinranking.com/coin/SbWqqTui-+energyswap-ens|~|19.09|~|  +0.56%  
51|~|BabyDogeCoin|~|https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge|~|0.0₆0237|~|  +1.39%  
"@

$lines = $rowsData -split "`n" | Where-Object { $_.Trim().Length -gt 0 }

foreach ($line in $lines) {
    $parts = $line -split "\|~\|"
    $r = [int]$parts[0]
    $coin = $parts[1]
    $link = $parts[2]
    $price = $parts[3]
    $volume = $parts[4]

    $ws.Range("B$r").Value = $coin
    $ws.Range("C$r").Value = $link

    # Price strings (e.g. "63.955.26", "1.00", "0.0590") must stay text - Excel's
    # COM Value setter otherwise parses them as numbers/dates and mangles them.
    $priceCell = $ws.Range("D$r")
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"

    $ws.Range("E$r").Value = $volume
}
